# Apply updated crypto market data (prices + 1h volume change) per the
# Thu Aug  8 04:56:22 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" cells are numeric-looking strings (e.g. "484.21", "1.00")
# that must stay stored as TEXT (matching the original inline-string cells)
# instead of being auto-coerced to numbers by COM, so force a Text format
# on those cells before writing their new value.
$textCells = @('D5','D6','D7','D8','D10','D12','D16','D19','D20','D21','D22','D24','D25','D26','D29','D32','D33','D34','D35','D37','D38','D39','D40','D41','D42','D44','D45','D46','D48','D49','D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row-by-row cell updates (Coin / Link / Price / Volume(1h))
$ws.Range('D2').Value = '56.575.17'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '2.404.04'
$ws.Range('E3').Value = '  -3.87%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '484.21'
$ws.Range('E5').Value = '  -2.41%  '
$ws.Range('D6').Value = '152.16'
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  +16.15%  '
$ws.Range('D9').Value = '2.421.37'
$ws.Range('E9').Value = '  -3.62%  '
$ws.Range('D10').Value = '0.0992'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('E11').Value = '  -2.88%  '
$ws.Range('D12').Value = '0.333'
$ws.Range('E12').Value = '  -1.03%  '
$ws.Range('E13').Value = '  +1.31%  '
$ws.Range('D14').Value = '2.829.71'
$ws.Range('E14').Value = '  -3.52%  '
$ws.Range('D15').Value = '56.859.47'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '20.66'
$ws.Range('E16').Value = '  -4.06%  '
$ws.Range('E17').Value = '  -3.38%  '
$ws.Range('D18').Value = '2.423.61'
$ws.Range('E18').Value = '  -3.44%  '
$ws.Range('D19').Value = '4.70'
$ws.Range('E19').Value = '  +3.25%  '
$ws.Range('D20').Value = '323.16'
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('D21').Value = '9.93'
$ws.Range('E21').Value = '  -4.28%  '
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('D24').Value = '57.80'
$ws.Range('E24').Value = '  -2.10%  '
$ws.Range('D25').Value = '0.406'
$ws.Range('E25').Value = '  -1.55%  '
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('E27').Value = '  -3.46%  '
$ws.Range('D28').Value = '2.519.42'
$ws.Range('E28').Value = '  -3.56%  '
$ws.Range('D29').Value = '7.23'
$ws.Range('E29').Value = '  -6.00%  '
$ws.Range('D30').Value = '0.0₃0779'
$ws.Range('E30').Value = '  -4.61%  '
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').Value = '18.56'
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('D33').Value = '148.62'
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('D34').Value = '1.51'
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('D35').Value = '5.36'
$ws.Range('E35').Value = '  +1.70%  '
$ws.Range('E36').Value = '  -2.86%  '
$ws.Range('D37').Value = '3.68'
$ws.Range('E37').Value = '  -3.54%  '
$ws.Range('D38').Value = '0.841'
$ws.Range('E38').Value = '  -4.31%  '
$ws.Range('D39').Value = '0.101'
$ws.Range('E39').Value = '  +9.00%  '
$ws.Range('D40').Value = '34.05'
$ws.Range('E40').Value = '  -0.80%  '
$ws.Range('D41').Value = '3.51'
$ws.Range('E41').Value = '  -0.55%  '
$ws.Range('D42').Value = '1.36'
$ws.Range('E42').Value = '  -2.33%  '
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('D44').Value = '0.589'
$ws.Range('E44').Value = '  -4.12%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '264.59'
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0529'
$ws.Range('E46').Value = '  -6.84%  '
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0227'
$ws.Range('E48').Value = '  -1.86%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '4.64'
$ws.Range('E49').Value = '  -6.39%  '
$ws.Range('D50').Value = '17.33'
$ws.Range('E50').Value = '  -3.33%  '
$ws.Range('D51').Value = '1.859.02'
$ws.Range('E51').Value = '  -2.64%  '
